# Apply "Three FNN results added" update to the ASR Results sheet.
# This updates several B-column words and C-column counts to new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value
$changes = @{
    "C2"  = 13
    "C3"  = 12
    "C4"  = 3
    "C5"  = 11
    "B6"  = "<day>"
    "C6"  = 11
    "B7"  = "<five>"
    "C7"  = 8
    "C8"  = 10
    "B9"  = "<do>"
    "C9"  = 13
    "B10" = "<then>"
    "C10" = 6
    "C11" = 14
    "C13" = 11
    "C14" = 8
    "C15" = 7
    "C16" = 9
    "C17" = 12
    "B18" = "<a>"
    "C18" = 12
    "C19" = 17
    "C20" = 10
    "C22" = 6
    "C23" = 10
    "C24" = 9
    "C26" = 5
    "C28" = 12
    "C29" = 8
    "B30" = "<to>"
    "C30" = 7
    "C31" = 14
    "B32" = "<make>"
    "C32" = 11
    "B33" = "<by>"
    "C33" = 9
    "B34" = "<part>"
    "C34" = 12
    "C35" = 5
    "B36" = "<its>"
    "C36" = 15
    "C37" = 13
    "C38" = 6
    "C39" = 8
    "C40" = 10
    "C41" = 7
    "C42" = 7
    "C43" = 9
    "C45" = 12
    "B46" = "<there>"
    "C46" = 14
    "C47" = 10
    "B48" = "<a>"
    "C48" = 16
    "C49" = 10
    "B50" = "<xoxtroand>"
    "C50" = 6
    "B52" = "<its>"
    "C52" = 6
}

foreach ($ref in $changes.Keys) {
    $ws.Range($ref).Value = $changes[$ref]
}
